$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 32 new rows at the bottom of the data block (rows 50-81) to make room
# for the extra historical risk groups, mirroring the values that used to sit
# in rows 6-37 (the 8 oldest groups at the time), then rewrite rows 2-49 with
# the refreshed/shifted monthly data.
$ws.Rows("50:81").Insert()

$data = @(
    @(2, 44958, "Heightened", 0),
    @(3, 44958, "Low", 2),
    @(4, 44958, "Medium", 0),
    @(5, 44958, "Standard", 0),
    @(6, 44986, "Heightened", 0),
    @(7, 44986, "Low", 0),
    @(8, 44986, "Medium", 0),
    @(9, 44986, "Standard", 0),
    @(10, 45017, "Heightened", 0),
    @(11, 45017, "Low", 0),
    @(12, 45017, "Medium", 0),
    @(13, 45017, "Standard", 0),
    @(14, 45047, "Heightened", 0),
    @(15, 45047, "Low", 0),
    @(16, 45047, "Medium", 0),
    @(17, 45047, "Standard", 0),
    @(18, 45078, "Heightened", 0),
    @(19, 45078, "Low", 0),
    @(20, 45078, "Medium", 0),
    @(21, 45078, "Standard", 0),
    @(22, 45108, "Heightened", 0),
    @(23, 45108, "Low", 0),
    @(24, 45108, "Medium", 0),
    @(25, 45108, "Standard", 0),
    @(26, 45139, "Heightened", 0),
    @(27, 45139, "Low", 0),
    @(28, 45139, "Medium", 0),
    @(29, 45139, "Standard", 0),
    @(30, 45170, "Heightened", 0),
    @(31, 45170, "Low", 0),
    @(32, 45170, "Medium", 0),
    @(33, 45170, "Standard", 0),
    @(34, 45200, "Heightened", 0),
    @(35, 45200, "Low", 0),
    @(36, 45200, "Medium", 0),
    @(37, 45200, "Standard", 0),
    @(38, 45231, "Heightened", 0),
    @(39, 45231, "Low", 0),
    @(40, 45231, "Medium", 0),
    @(41, 45231, "Standard", 0),
    @(42, 45261, "Heightened", 0),
    @(43, 45261, "Low", 0),
    @(44, 45261, "Medium", 0),
    @(45, 45261, "Standard", 0),
    @(46, 45292, "Heightened", 0),
    @(47, 45292, "Low", 0),
    @(48, 45292, "Medium", 0),
    @(49, 45292, "Standard", 0),
    @(50, 45323, "Heightened", 39),
    @(51, 45323, "Low", 181),
    @(52, 45323, "Medium", 0),
    @(53, 45323, "Standard", 33),
    @(54, 45352, "Heightened", 1),
    @(55, 45352, "Low", 2),
    @(56, 45352, "Medium", 0),
    @(57, 45352, "Standard", 0),
    @(58, 45383, "Heightened", 8),
    @(59, 45383, "Low", 140),
    @(60, 45383, "Medium", 0),
    @(61, 45383, "Standard", 9),
    @(62, 45413, "Heightened", 0),
    @(63, 45413, "Low", 8),
    @(64, 45413, "Medium", 1),
    @(65, 45413, "Standard", 0),
    @(66, 45444, "Heightened", 0),
    @(67, 45444, "Low", 83),
    @(68, 45444, "Medium", 0),
    @(69, 45444, "Standard", 0),
    @(70, 45474, "Heightened", 0),
    @(71, 45474, "Low", 1),
    @(72, 45474, "Medium", 0),
    @(73, 45474, "Standard", 0),
    @(74, 45505, "Heightened", 0),
    @(75, 45505, "Low", 1),
    @(76, 45505, "Medium", 0),
    @(77, 45505, "Standard", 0),
    @(78, 45536, "Heightened", 0),
    @(79, 45536, "Low", 1),
    @(80, 45536, "Medium", 0),
    @(81, 45536, "Standard", 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
}
